# Removed the "no stagger" options from the fleet reservation parameter sheet.
#
# Net effect (see commit message / diff):
#  - Rows 3 & 4 ("Base SLEP case with no stagger") keep their numbers but lose
#    their Notes label and have Stagger (col R) bumped from 0 to 0.5.
#  - Row 5 is repurposed into a new "Base case no SLEP with surge" scenario
#    (SLEP_or_not -> FALSE, Stagger -> 0, addHours -> 10800, TTR formula ->
#    24*30*12, surgetime formula -> 24*365*12, Notes -> new label).
#  - The old row 6 is removed outright, shifting old rows 7 & 8 up to become
#    the new rows 6 & 7. The sheet shrinks from 8 data rows to 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: bump Stagger to 0.5 and drop the "no stagger" note ---
$ws.Range("R3").Value = 0.5
$ws.Range("X3").ClearContents()

# --- Row 4: bump Stagger to 0.5 and drop the "no stagger" note ---
$ws.Range("R4").Value = 0.5
$ws.Range("X4").ClearContents()

# --- Row 5: turn into the new "Base case no SLEP with surge" scenario ---
$ws.Range("P5").Value = $false
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 10800
$ws.Range("T5").Formula = "=24*30*12"
$ws.Range("W5").Formula = "=24*365*12"
$ws.Range("X5").Value = "Base case no SLEP with surge"

# --- Remove the old row 6 entirely; rows 7 & 8 shift up to 6 & 7 ---
$ws.Rows.Item(6).Delete()
